$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A56").Value = "Rewrite board basic functions using validation"
$ws.Range("A56").Font.Bold = $true
$ws.Range("A57").Value = "lila\src\main\scala\model\Board.scala"
$ws.Range("B57").Value = "x"

$ws.Range("C57").Select() | Out-Null
